$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.592.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.699.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4035"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.550"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +8.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.83"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +15.13%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08800"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.310"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +11.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.634"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.696.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07027"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.935"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.583.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.977"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.338"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.40"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.81"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.235"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.631"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +19.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.113"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.884.60"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.496"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +13.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08571"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.984"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.17"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2760"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.78"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02776"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09036"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7756"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7297"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.52"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.510"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.190"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.300"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +13.80%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.32%  "
